$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated cell values
$ws.Range('D2').Value = '25.998.72'
$ws.Range('D3').Value = '1.741.40'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9997'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '247.31'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.44%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5052'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -4.31%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2750'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06186'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('B10').Value = 'TRON'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07269'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.42%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '1.743.08'
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.6546'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.92%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '15.10'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.681'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '77.59'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.9999'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').Value = '26.019.88'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.91'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('D21').Value = '1.973.52'
$ws.Range('E21').Value = '  +0.57%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.498'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.45%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.718'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.396'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.61%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '135.25'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.34%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.504'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '105.34'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.962'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.27%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08177'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.51%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.686'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.23%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04685'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.93%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9954'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.29%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6136'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.87%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.755'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.98%  '
$ws.Range('E38').Value = '  +1.26%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.929'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '100.94'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.33%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.3918'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.90%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.7641'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.86%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.012'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.30%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1164'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.81%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '6.310'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.43%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '55.55'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.44%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '30.74'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.630'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.80%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.3473'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.84%  '
